$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Small "cos(theta)" helper table that sits to the right of the existing
# measurement table (rows 10-16, columns B:F).
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "úhel [deg]"
$ws.Range("C10").Value = "cos"
$ws.Range("D10").Value = "úhel pm [deg]"
$ws.Range("T10").Font.Bold = $true

$ws.Range("B11").Value = 102
$ws.Range("C11").Formula = "=COS(RADIANS(B11))"
$ws.Range("D11").Value = 4
$ws.Range("E11").Formula = "=D11/B11"
$ws.Range("F11").Formula = "=E11*C11"

$ws.Range("B12").Value = 88
$ws.Range("C12").Formula = "=COS(RADIANS(B12))"
$ws.Range("D12").Value = 4
$ws.Range("E12").Formula = "=D12/B12"
$ws.Range("F12").Formula = "=E12*C12"

$ws.Range("B13").Value = 77
$ws.Range("C13").Formula = "=COS(RADIANS(B13))"
$ws.Range("D13").Value = 4
$ws.Range("E13").Formula = "=D13/B13"
$ws.Range("F13").Formula = "=E13*C13"

$ws.Range("B14").Value = 99
$ws.Range("C14").Formula = "=COS(RADIANS(B14))"
$ws.Range("D14").Value = 6
$ws.Range("E14").Formula = "=D14/B14"
$ws.Range("F14").Formula = "=E14*C14"

$ws.Range("B15").Value = 91
$ws.Range("C15").Formula = "=COS(RADIANS(B15))"
$ws.Range("D15").Value = 2
$ws.Range("E15").Formula = "=D15/B15"
$ws.Range("F15").Formula = "=E15*C15"

$ws.Range("B16").Value = 68
$ws.Range("C16").Formula = "=COS(RADIANS(B16))"
$ws.Range("D16").Value = 3
$ws.Range("E16").Formula = "=D16/B16"
$ws.Range("F16").Formula = "=E16*C16"

# ---------------------------------------------------------------------------
# First small "surface tension of test liquids" table (rows 18-24).
# Column A is filled top to bottom first (this drives the shared-string
# insertion order), then the header B18 label is added.
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "kapalina"
$ws.Range("A19").Value = "destilovaná voda"
$ws.Range("A20").Value = "etylenglykol"
$ws.Range("A21").Value = "dijodometan"
$ws.Range("A22").Value = "glycerol"
$ws.Range("A23").Value = "formamid"
$ws.Range("A24").Value = "alpha-bromnaftalen"

$ws.Range("B18").Value = "\gamma_l"
$ws.Range("B19").Value = 72.8
$ws.Range("B20").Value = 47.7
$ws.Range("B21").Value = 50.8
$ws.Range("B22").Value = 64
$ws.Range("B23").Value = 58.2
$ws.Range("B24").Value = 44.4

$ws.Range("D18").Font.Bold = $true

$ws.Range("A18").Font.Bold = $true
$ws.Range("B18").Font.Bold = $true

# ---------------------------------------------------------------------------
# Second table (rows 26-32): contact-angle measurements + derived cos(theta)
# for each test liquid. Column F (pm cos theta) is written before column E
# (relative pm) so the new shared strings land on the same indices as the
# original edit.
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = "kapalina"
$ws.Range("B26").Value = "\theta [deg]"
$ws.Range("C26").Value = "pm \theta [deg]"
$ws.Range("D26").Value = "cos \theta"
$ws.Range("F26").Value = "pm cos \theta"
$ws.Range("E26").Value = "relative pm"

$ws.Range("A26").Font.Bold = $true
$ws.Range("B26").Font.Bold = $true
$ws.Range("C26").Font.Bold = $true
$ws.Range("D26").Font.Bold = $true
$ws.Range("E26").Font.Bold = $true
$ws.Range("F26").Font.Bold = $true

$ws.Range("A27").Value = "destilovaná voda"
$ws.Range("B27").Value = 102.4
$ws.Range("C27").Value = 4
$ws.Range("D27").Formula = "=COS(RADIANS(B27))"
$ws.Range("E27").Formula = "=C27/B27"
$ws.Range("F27").Formula = "=E27*D27"

$ws.Range("A28").Value = "etylenglykol"
$ws.Range("B28").Value = 88.1
$ws.Range("C28").Value = 4
$ws.Range("D28").Formula = "=COS(RADIANS(B28))"
$ws.Range("E28").Formula = "=C28/B28"
$ws.Range("F28").Formula = "=E28*D28"

$ws.Range("A29").Value = "dijodometan"
$ws.Range("B29").Value = 76.7
$ws.Range("C29").Value = 4
$ws.Range("D29").Formula = "=COS(RADIANS(B29))"
$ws.Range("E29").Formula = "=C29/B29"
$ws.Range("F29").Formula = "=E29*D29"

$ws.Range("A30").Value = "glycerol"
$ws.Range("B30").Value = 99.5
$ws.Range("C30").Value = 6
$ws.Range("D30").Formula = "=COS(RADIANS(B30))"
$ws.Range("E30").Formula = "=C30/B30"
$ws.Range("F30").Formula = "=E30*D30"

$ws.Range("A31").Value = "formamid"
$ws.Range("B31").Value = 90.5
$ws.Range("C31").Value = 2
$ws.Range("D31").Formula = "=COS(RADIANS(B31))"
$ws.Range("E31").Formula = "=C31/B31"
$ws.Range("F31").Formula = "=E31*D31"

$ws.Range("A32").Value = "alpha-bromnaftalen"
$ws.Range("B32").Value = 67.7
$ws.Range("C32").Value = 3
$ws.Range("D32").Formula = "=COS(RADIANS(B32))"
$ws.Range("E32").Formula = "=C32/B32"
$ws.Range("F32").Formula = "=E32*D32"

# ---------------------------------------------------------------------------
# Cosmetics: column widths for the newly populated columns, and move the
# active selection the way the author left it.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 19.33203125
$ws.Columns("C").ColumnWidth = 14.5546875
$ws.Columns("D").ColumnWidth = 12.6640625
$ws.Columns("E").ColumnWidth = 12.77734375
$ws.Columns("F").ColumnWidth = 12.77734375

$ws.Range("C17").Select()

Write-Host "Done"
